$wb = $excel.ActiveWorkbook

# --- Sheet ALC: 45 cell change(s) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value2 = 4167
$ws.Range("I28").Value2 = 1121.6
$ws.Range("J28").Value2 = 6070.375
$ws.Range("K28").Value2 = 1121.6
$ws.Range("L28").Value2 = 6070.375
$ws.Range("M28").Value2 = -636.5999999999999
$ws.Range("N28").Value2 = -7040.375
$ws.Range("H40").Value2 = 4794.2046
$ws.Range("I40").Value2 = 3977.2727
$ws.Range("K40").Value2 = 3977.2727
$ws.Range("M40").Value2 = -3802.2727
$ws.Range("H64").Value2 = 10000
$ws.Range("I64").Value2 = 10000
$ws.Range("J64").Value2 = 10000
$ws.Range("K64").Value2 = 10000
$ws.Range("L64").Value2 = 10000
$ws.Range("M64").Value2 = -9752
$ws.Range("N64").Value2 = -10496
$ws.Range("H67").Value2 = 10000
$ws.Range("I67").Value2 = 10000
$ws.Range("J67").Value2 = 10000
$ws.Range("K67").Value2 = 10000
$ws.Range("L67").Value2 = 10000
$ws.Range("M67").Value2 = -9142
$ws.Range("N67").Value2 = -11716
$ws.Range("H74").Value2 = 3498.3333
$ws.Range("I74").Value2 = 3498.3333
$ws.Range("K74").Value2 = 3498.3333
$ws.Range("M74").Value2 = -2562.3333
$ws.Range("H77").Value2 = 3498.3333
$ws.Range("I77").Value2 = 3498.3333
$ws.Range("K77").Value2 = 17491.6665
$ws.Range("M77").Value2 = -12811.6665
$ws.Range("H103").Value2 = 2819.25
$ws.Range("I103").Value2 = 2683.1
$ws.Range("K103").Value2 = 8049.299999999999
$ws.Range("M103").Value2 = -7463.299999999999
$ws.Range("H107").Value2 = 1043.25
$ws.Range("I107").Value2 = 1038.1111
$ws.Range("K107").Value2 = 1038.1111
$ws.Range("M107").Value2 = 881.8888999999999
$ws.Range("H132").Value2 = 9509.357
$ws.Range("J132").Value2 = 13590.333
$ws.Range("L132").Value2 = 40770.999
$ws.Range("N132").Value2 = -45830.999

# --- Sheet ARM: 70 cell change(s) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value2 = 7665.4287
$ws.Range("I16").Value2 = 12814.5
$ws.Range("K16").Value2 = 12814.5
$ws.Range("M16").Value2 = -12527.5
$ws.Range("H32").Value2 = 4727.65
$ws.Range("I32").Value2 = 4713.316
$ws.Range("K32").Value2 = 4713.316
$ws.Range("M32").Value2 = -4426.316
$ws.Range("H45").Value2 = 2628
$ws.Range("I45").Value2 = 1856
$ws.Range("K45").Value2 = 1856
$ws.Range("M45").Value2 = -1479
$ws.Range("H61").Value2 = 6486.5625
$ws.Range("I61").Value2 = 5473.125
$ws.Range("J61").Value2 = 7500
$ws.Range("K61").Value2 = 5473.125
$ws.Range("L61").Value2 = 7500
$ws.Range("M61").Value2 = -5261.125
$ws.Range("N61").Value2 = -7924
$ws.Range("H63").Value2 = 11474.125
$ws.Range("J63").Value2 = 20924.75
$ws.Range("L63").Value2 = 20924.75
$ws.Range("N63").Value2 = -22296.75
$ws.Range("H66").Value2 = 11474.125
$ws.Range("J66").Value2 = 20924.75
$ws.Range("L66").Value2 = 104623.75
$ws.Range("N66").Value2 = -111487.75
$ws.Range("H74").Value2 = 2646.077
$ws.Range("I74").Value2 = 2646.077
$ws.Range("J74").Value2 = 0
$ws.Range("K74").Value2 = 2646.077
$ws.Range("L74").Value2 = 0
$ws.Range("M74").Value2 = -1772.077
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value2 = 2646.077
$ws.Range("I77").Value2 = 2646.077
$ws.Range("J77").Value2 = 0
$ws.Range("K77").Value2 = 13230.385
$ws.Range("L77").Value2 = 0
$ws.Range("M77").Value2 = -8862.385000000002
$ws.Range("N77").ClearContents()
$ws.Range("H97").Value2 = 515.53845
$ws.Range("I97").Value2 = 643.2
$ws.Range("J97").Value2 = 90
$ws.Range("K97").Value2 = 643.2
$ws.Range("L97").Value2 = 90
$ws.Range("M97").Value2 = -147.2
$ws.Range("N97").Value2 = -1082
$ws.Range("H102").Value2 = 4366.5
$ws.Range("I102").Value2 = 3225.889
$ws.Range("J102").Value2 = 6419.6
$ws.Range("K102").Value2 = 3225.889
$ws.Range("L102").Value2 = 6419.6
$ws.Range("M102").Value2 = -1603.889
$ws.Range("N102").Value2 = -9663.6
$ws.Range("H110").Value2 = 932.3125
$ws.Range("I110").Value2 = 685.3077
$ws.Range("K110").Value2 = 685.3077
$ws.Range("M110").Value2 = 1359.6923
$ws.Range("H132").Value2 = 2606.6667
$ws.Range("I132").Value2 = 2451.96
$ws.Range("K132").Value2 = 7355.88
$ws.Range("M132").Value2 = -4825.88
$ws.Range("H136").Value2 = 6486.5625
$ws.Range("I136").Value2 = 5473.125
$ws.Range("J136").Value2 = 7500
$ws.Range("K136").Value2 = 16419.375
$ws.Range("L136").Value2 = 22500
$ws.Range("M136").Value2 = -13869.375
$ws.Range("N136").Value2 = -27600

# --- Sheet BSM: 22 cell change(s) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value2 = 5667.4
$ws.Range("I86").Value2 = 4686.6
$ws.Range("J86").Value2 = 6648.2
$ws.Range("K86").Value2 = 4686.6
$ws.Range("L86").Value2 = 6648.2
$ws.Range("M86").Value2 = -3563.6
$ws.Range("N86").Value2 = -8894.200000000001
$ws.Range("H89").Value2 = 5667.4
$ws.Range("I89").Value2 = 4686.6
$ws.Range("J89").Value2 = 6648.2
$ws.Range("K89").Value2 = 23433
$ws.Range("L89").Value2 = 33241
$ws.Range("M89").Value2 = -17817
$ws.Range("N89").Value2 = -44473
$ws.Range("H94").Value2 = 1866.3334
$ws.Range("I94").Value2 = 1866.3334
$ws.Range("K94").Value2 = 1866.3334
$ws.Range("M94").Value2 = -1415.3334
$ws.Range("H134").Value2 = 3227.182
$ws.Range("I134").Value2 = 3227.182
$ws.Range("K134").Value2 = 9681.545999999998
$ws.Range("M134").Value2 = -7146.545999999998

# --- Sheet CRP: 23 cell change(s) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value2 = 3666.3572
$ws.Range("I58").Value2 = 3115.3333
$ws.Range("K58").Value2 = 3115.3333
$ws.Range("M58").Value2 = -2912.3333
$ws.Range("H62").Value2 = 4647.143
$ws.Range("I62").Value2 = 4647.143
$ws.Range("K62").Value2 = 4647.143
$ws.Range("M62").Value2 = -4023.143
$ws.Range("H65").Value2 = 4647.143
$ws.Range("I65").Value2 = 4647.143
$ws.Range("K65").Value2 = 23235.715
$ws.Range("M65").Value2 = -20115.715
$ws.Range("H105").Value2 = 2998.8
$ws.Range("I105").Value2 = 2997.25
$ws.Range("J105").Value2 = 2999.8333
$ws.Range("K105").Value2 = 2997.25
$ws.Range("L105").Value2 = 2999.8333
$ws.Range("M105").Value2 = -1250.25
$ws.Range("N105").Value2 = -6493.8333
$ws.Range("H136").Value2 = 3666.3572
$ws.Range("I136").Value2 = 3115.3333
$ws.Range("K136").Value2 = 9345.999899999999
$ws.Range("M136").Value2 = -6795.999899999999

# --- Sheet CUL: 29 cell change(s) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value2 = 0
$ws.Range("I11").Value2 = 0
$ws.Range("J11").Value2 = 0
$ws.Range("K11").Value2 = 0
$ws.Range("L11").Value2 = 0
$ws.Range("M11").ClearContents()
$ws.Range("N11").ClearContents()
$ws.Range("H26").Value2 = 44
$ws.Range("I26").Value2 = 60
$ws.Range("J26").Value2 = 33.333332
$ws.Range("K26").Value2 = 180
$ws.Range("L26").Value2 = 99.999996
$ws.Range("M26").Value2 = 108
$ws.Range("N26").Value2 = -675.999996
$ws.Range("H37").Value2 = 45000
$ws.Range("J37").Value2 = 45000
$ws.Range("L37").Value2 = 135000
$ws.Range("N37").Value2 = -135224
$ws.Range("H92").Value2 = 760.5
$ws.Range("I92").Value2 = 662.3333
$ws.Range("J92").Value2 = 1055
$ws.Range("K92").Value2 = 1986.9999
$ws.Range("L92").Value2 = 3165
$ws.Range("M92").Value2 = -738.9999
$ws.Range("N92").Value2 = -5661
$ws.Range("H104").Value2 = 9402.308000000001
$ws.Range("I104").Value2 = 7375
$ws.Range("K104").Value2 = 22125
$ws.Range("M104").Value2 = -19504

# --- Sheet GSM: 12 cell change(s) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value2 = 1489.6666
$ws.Range("I70").Value2 = 1489.6666
$ws.Range("K70").Value2 = 1489.6666
$ws.Range("M70").Value2 = -1219.6666
$ws.Range("H73").Value2 = 1489.6666
$ws.Range("I73").Value2 = 1489.6666
$ws.Range("K73").Value2 = 1489.6666
$ws.Range("M73").Value2 = -553.6666
$ws.Range("H132").Value2 = 4215.75
$ws.Range("I132").Value2 = 2354.3333
$ws.Range("K132").Value2 = 7062.999899999999
$ws.Range("M132").Value2 = -4532.999899999999

# --- Sheet LTW: 19 cell change(s) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value2 = 765.7143
$ws.Range("I55").Value2 = 814.5833
$ws.Range("J55").Value2 = 472.5
$ws.Range("K55").Value2 = 814.5833
$ws.Range("L55").Value2 = 472.5
$ws.Range("M55").Value2 = -641.5833
$ws.Range("N55").Value2 = -818.5
$ws.Range("H61").Value2 = 3281.3333
$ws.Range("I61").Value2 = 922
$ws.Range("K61").Value2 = 922
$ws.Range("M61").Value2 = -720
$ws.Range("H113").Value2 = 3281.3333
$ws.Range("I113").Value2 = 922
$ws.Range("K113").Value2 = 922
$ws.Range("M113").Value2 = 1248
$ws.Range("H133").Value2 = 100000
$ws.Range("J133").Value2 = 100000
$ws.Range("L133").Value2 = 100000
$ws.Range("N133").Value2 = -105060

# --- Sheet WVR: 15 cell change(s) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value2 = 19523
$ws.Range("J105").Value2 = 19523
$ws.Range("L105").Value2 = 19523
$ws.Range("N105").Value2 = -26511
$ws.Range("H113").Value2 = 638.375
$ws.Range("I113").Value2 = 641.6
$ws.Range("K113").Value2 = 1924.8
$ws.Range("M113").Value2 = 245.1999999999998
$ws.Range("H132").Value2 = 1931.4
$ws.Range("I132").Value2 = 1931.4
$ws.Range("J132").Value2 = 0
$ws.Range("K132").Value2 = 5794.200000000001
$ws.Range("L132").Value2 = 0
$ws.Range("M132").Value2 = -3264.200000000001
$ws.Range("N132").ClearContents()
